# Horizontal navigation bar added
# - Update the "Date" column (column B) entries that read 01/08/2023 /
#   2023/08/01 / 08/01/2023 to the 3rd of August instead of the 1st
#   (03/08/2023 / 2023/08/03 / 08/03/2023), preserving each cell's existing
#   text format (quote-prefixed / literal text).
# - Normalize the format of the last two rows (B16/B17) to match the plain
#   text style used by the rest of the column.
# - Move the active selection down to the bottom of the sheet (B17) and
#   scroll the view down so row 4 is at the top (acting as a simple
#   "navigation" jump to the end of the data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that currently display 01/08/2023 -> 03/08/2023
$augThirdDMY = @("B1", "B2", "B3", "B4", "B5", "B6", "B7", "B8", "B10", "B11", "B16", "B17")
foreach ($addr in $augThirdDMY) {
    $ws.Range($addr).Value = "'03/08/2023"
}

# Cell that currently displays 2023/08/01 -> 2023/08/03
$ws.Range("B12").Value = "'2023/08/03"

# Cell that currently displays 08/01/2023 -> 08/03/2023
$ws.Range("B14").Value = "'08/03/2023"

# B16 and B17 were still using the quote-prefixed date style (s=9); align
# them with the plain quote-prefixed text style used elsewhere in column B.
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B17").NumberFormat = "@"

# Re-apply the values now that the number format matches, so the stored
# text / style stay consistent.
$ws.Range("B16").Value = "'03/08/2023"
$ws.Range("B17").Value = "'03/08/2023"

# Navigate: scroll so row 4 is at the top and select the last cell (B17),
# mimicking a "jump to bottom" navigation action.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("B17").Select()
